# Insert a new data row at row 21 (pushing existing rows 21-43 down to 22-44)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 21:43 down to 22:44, leaving a blank row 21 ready to fill in.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new "Feria Lagunitas de Puerto
# Montt" / "Haba" record. Columns that are constant for every record in this
# sheet (A, B, C, E, F, G, H, I, N, Q, R) are copied from the surrounding rows.
$ws.Range("A21").Value = 4
$ws.Range("B21").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C21").Value = "Los Lagos"
$ws.Range("D21").Value = 44484
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 100112026
$ws.Range("G21").Value = "Haba"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = 10000
$ws.Range("N21").Value = "$/saco 25 kilos"
$ws.Range("O21").Value = "Región Metropolitana"
$ws.Range("P21").Value = 400
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
